# Auto-generated edit script for R_AHP_CasoBase.xlsx
# Applies updated AHP weights/rankings per the commit diff.

$wb = $excel.ActiveWorkbook

# --- Pesos_Locales_Económico : updated local weights (column B) ---
$ws = $wb.Worksheets.Item("Pesos_Locales_Económico")
$ws.Range("B2").Value = 0.06859393436079969
$ws.Range("B3").Value = 0.06859393436079969
$ws.Range("B4").Value = 0.1451001396860344
$ws.Range("B5").Value = 0.06859393436079965
$ws.Range("B6").Value = 0.1451001396860344
$ws.Range("B7").Value = 0.1451001396860344
$ws.Range("B8").Value = 0.01428744057464852
$ws.Range("B9").Value = 0.01428744057464852
$ws.Range("B10").Value = 0.06859393436079965
$ws.Range("B11").Value = 0.0697866565851307
$ws.Range("B12").Value = 0.008867907233478526
$ws.Range("B13").Value = 0.1451001396860344
$ws.Range("B14").Value = 0.009419377695460395
$ws.Range("B15").Value = 0.01428744057464852
$ws.Range("B16").Value = 0.01428744057464852

# --- Ranking_Alternativas : re-sorted alternative order + updated global weights ---
$ws = $wb.Worksheets.Item("Ranking_Alternativas")
$ws.Range("A2").Value = "Placilla"
$ws.Range("B2").Value = 0.1061177752378623
$ws.Range("A3").Value = "Plaza Justicia"
$ws.Range("B3").Value = 0.1036947702478947
$ws.Range("A4").Value = "Jean y Marie Thierry"
$ws.Range("B4").Value = 0.0864644866674637
$ws.Range("A5").Value = "Marcelo Mena"
$ws.Range("B5").Value = 0.08513647894959353
$ws.Range("B6").Value = 0.06602550478671811
$ws.Range("A7").Value = "Quebrada Verde"
$ws.Range("B7").Value = 0.06590218896974658
$ws.Range("A8").Value = "Placeres"
$ws.Range("B8").Value = 0.06578654256230929
$ws.Range("A9").Value = "Laguna Verde"
$ws.Range("B9").Value = 0.06396582046080838
$ws.Range("A10").Value = "Puertas Negras"
$ws.Range("B10").Value = 0.06374873672567684
$ws.Range("B11").Value = 0.05960176597705624
$ws.Range("A12").Value = "Reina Isabel 2"
$ws.Range("B12").Value = 0.05574321597212605
$ws.Range("A13").Value = "Esperanza"
$ws.Range("B13").Value = 0.05267876917419514
$ws.Range("B14").Value = 0.05217500779969719
$ws.Range("B15").Value = 0.04258198851215252
$ws.Range("B16").Value = 0.03037694795669948

# --- Resultados : updated global weights (column B) ---
$ws = $wb.Worksheets.Item("Resultados")
$ws.Range("B2").Value = 0.05217500779969719
$ws.Range("B3").Value = 0.05960176597705624
$ws.Range("B4").Value = 0.05267876917419514
$ws.Range("B5").Value = 0.0864644866674637
$ws.Range("B6").Value = 0.06396582046080838
$ws.Range("B7").Value = 0.06602550478671811
$ws.Range("B8").Value = 0.08513647894959353
$ws.Range("B9").Value = 0.03037694795669948
$ws.Range("B10").Value = 0.06578654256230929
$ws.Range("B11").Value = 0.1061177752378623
$ws.Range("B12").Value = 0.1036947702478947
$ws.Range("B13").Value = 0.06374873672567684
$ws.Range("B14").Value = 0.06590218896974658
$ws.Range("B15").Value = 0.05574321597212605
$ws.Range("B16").Value = 0.04258198851215252

# --- Matriz_Económico : updated pairwise comparison entries ---
$ws = $wb.Worksheets.Item("Matriz_Económico")
$ws.Range("D2").Value = 0.3333333333333333
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.3333333333333333
$ws.Range("N2").Value = 7
$ws.Range("D3").Value = 0.3333333333333333
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.3333333333333333
$ws.Range("N3").Value = 7
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 3
$ws.Range("E4").Value = 3
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 3
$ws.Range("D5").Value = 0.3333333333333333
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.3333333333333333
$ws.Range("N5").Value = 7
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 3
$ws.Range("E6").Value = 3
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 3
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 3
$ws.Range("E7").Value = 3
$ws.Range("J7").Value = 3
$ws.Range("K7").Value = 3
$ws.Range("D10").Value = 0.3333333333333333
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.3333333333333333
$ws.Range("N10").Value = 7
$ws.Range("D11").Value = 0.3333333333333333
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.3333333333333333
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = 3
$ws.Range("E13").Value = 3
$ws.Range("J13").Value = 3
$ws.Range("K13").Value = 3
$ws.Range("B14").Value = 0.1428571428571428
$ws.Range("C14").Value = 0.1428571428571428
$ws.Range("E14").Value = 0.1428571428571428
$ws.Range("J14").Value = 0.1428571428571428

